$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 5: code_id field - change type from integer to string, update hint and constraint regex
$ws.Range("A5").Value = "string"
$ws.Range("D5").Value = 'Exemple "190-123456"'
$ws.Range("F5").Value = "regex(., '^[0-9]{3}-[0-9]{6}$')"

# Row 6: code_id2 field - change type from integer to string, update hint
$ws.Range("A6").Value = "string"
$ws.Range("D6").Value = 'Exemple "190-123456"'
